$wb = $excel.ActiveWorkbook

# Switch to the Materials sheet (this becomes the active/selected tab)
$ws = $wb.Worksheets.Item("Materials")
$ws.Activate()

# Rename "Spaghetti Box" -> "Spaghetti Boxes" and bump its Max Checkout Quantity 2 -> 3
$ws.Range("A4").Value = "Spaghetti Boxes"
$ws.Range("C4").Value = 3

# Widen column A so the longer item names still fit (matches bestFit behavior)
$ws.Columns.Item(1).ColumnWidth = 15

# Add a new material row: Marshmellow Bag
[void]$ws.Range("A7").Select()
$ws.Cells.Item(7, 1).Value = "Marshmellow Bag"
$ws.Cells.Item(7, 2).Value = 3
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0

# Leave selection on the row after the new entry, like typing down through the row
[void]$ws.Range("A8").Select()
